$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Warnings")
$ws2 = $wb.Worksheets.Item("Validation")

$ws1.Range('C13').Value2 = 'Result*'
$ws1.Range('B14').Value2 = 'Step1*'
$ws1.Range('B15').Value2 = 'Step2*'
$ws1.Range('B16').Value2 = 'Step3*'
$ws1.Range('B17').Value2 = 'Step4*'
$ws1.Range('B18').Value2 = 'Step5*'
$ws1.Range('C27').Value2 = 'Step1*'
$ws1.Range('D27').Value2 = 'Step2*'
$ws1.Range('E27').Value2 = 'Step3*'
$ws1.Range('F27').Value2 = 'Step4*'
$ws1.Range('G27').Value2 = 'Step5*'
$ws1.Range('D37').Value2 = 'OtherValues*'
$ws1.Range('B44').Value2 = 'Step1*'
$ws1.Range('B46').Value2 = 'Step3*'
$ws1.Range('C49').Value2 = 'Result0*'
$ws1.Range('D49').Value2 = 'Result*'
$ws1.Range('B50').Value2 = 'Step1*'
$ws1.Range('B51').Value2 = 'Step2*'
$ws1.Range('B52').Value2 = 'Step3*'
$ws1.Range('B53').Value2 = 'Step4*'
$ws1.Range('B54').Value2 = 'Step5*'
$ws1.Range('C62').Value2 = 'Res_ult'
$ws1.Range('D62').Value2 = 'Res'
$ws1.Range('B63').Value2 = 'ult_Value'
$ws1.Range('B64').Value2 = 'Value'
$ws2.Range('C6').Value2 = 'Values*'
$ws2.Range('B7').Value2 = '1Step1*'
$ws2.Range('B8').Value2 = 'My Step*'
$ws2.Range('B9').Value2 = ' Step_1 *'
$ws2.Range('B10').Value2 = 'Step1* *'
$ws2.Range('B11').Value2 = 'Step^*'
$ws2.Range('B12').Value2 = 'Step-1*'
$ws2.Range('B13').Value2 = 'Step.1*'
$ws2.Range('B14').Value2 = 'Step__1*'
$ws2.Range('C19').Value2 = 'Values*'
$ws2.Range('B20').Value2 = 'тест*'
$ws2.Range('B21').Value2 = 'アスタリスク*'
$ws2.Range('B22').Value2 = 'TestЫЫЫ*'
$ws2.Range('B23').Value2 = ' 印刷中的星号*'
$ws2.Range('B24').Value2 = 'aaa*'
$ws2.Range('B30').Value2 = 'myValue*'
$ws2.Range('B31').Value2 = 'MyValue*'
$ws2.Range('B32').Value2 = 'result*'
